# Update countries & provincias Spain
# - Re-sort several tied/changed country rows in the "Pais" ranking table
# - Update a few case-count figures (Australia, Butan/Gambia/Antigua y Barbuda)
# - Bump the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refresh timestamp string
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 05:05"

# Row 62: Australia - update active/recovered counts
$ws.Cells.Item(62, 4).Value = 6531
$ws.Cells.Item(62, 5).Value = 481

# Rows 189-191: Butan overtakes Gambia and Antigua y Barbuda in total cases
$ws.Cells.Item(189, 1).Value = "Butan"
$ws.Cells.Item(189, 2).Value = 27
$ws.Cells.Item(189, 3).Value = 3
$ws.Cells.Item(189, 4).Value = 6
$ws.Cells.Item(189, 5).Value = 21
$ws.Cells.Item(189, 8).Value = 0

$ws.Cells.Item(190, 1).Value = "Gambia"
$ws.Cells.Item(190, 2).Value = 25
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 13
$ws.Cells.Item(190, 5).Value = 11
$ws.Cells.Item(190, 8).Value = 1

$ws.Cells.Item(191, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(191, 2).Value = 25
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 19
$ws.Cells.Item(191, 5).Value = 3
$ws.Cells.Item(191, 8).Value = 3

# Rows 198-200: Santa Lucia and Nueva Caledonia swap (tied figures)
$ws.Cells.Item(198, 1).Value = "Santa Lucia"
$ws.Cells.Item(200, 1).Value = "Nueva Caledonia"

# Rows 215-216: San Bartolome and Bonaire, San Eustaquio y Saba swap (tied figures)
$ws.Cells.Item(215, 1).Value = "San Bartolome"
$ws.Cells.Item(216, 1).Value = "Bonaire, San Eustaquio y Saba"
